# Update the "Conditional R2" and "Marginal R2" values for the two
# lme_disp_di rows of the results table.
#
# Table layout (1-indexed):
#   Row 1: header (Model, Parameter, Estimate, Std.Error, Conditional R2,
#                  Marginal R2, DF, t-value, p-value, n, K, Res.LL, AICc,
#                  ΔAICc, AIC weight)
#   Row 2: lme_disp_di / (Intercept)               -> Conditional R2 = 0.41, Marginal R2 = 0
#   Row 3: lme_disp_di / DispersalPotentialKmY      -> Conditional R2 = 0.41, Marginal R2 = 0
#   Row 4: lme_cv_di   / (Intercept)
#   Row 5: lme_cv_di   / q3ClimVeloKmY_RelScale
#
# Conditional R2 is column 5, Marginal R2 is column 6.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

foreach ($rowIndex in 2, 3) {
    $condCell = $tbl.Cell($rowIndex, 5)
    $condRange = $condCell.Range
    $condRange.End = $condRange.End - 1
    if ($condRange.Text -eq "0.41") {
        $condRange.Text = "0.4"
    }

    $margCell = $tbl.Cell($rowIndex, 6)
    $margRange = $margCell.Range
    $margRange.End = $margRange.End - 1
    if ($margRange.Text -eq "0") {
        $margRange.Text = "0.02"
    }
}
